$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated TPM-derived values for rows 2-26 (per the "update scripts wuth new tpm" commit).
# Each entry maps a row number to the column/value pairs that changed for that row.
$updates = @{
    2 = @{ "G" = 8.542726333333333; "H" = 25.628179; "I" = 0.2146499313812649; "J" = 0.2146499313812649; "M" = 101.1420973333333; "N" = 303.426292; "O" = 0.6720924517779291; "P" = 0.6720924517779291; "Q" = 864.0292582980296; "R" = 7776.263324682266; "S" = 0.1442645986559986; "T" = 0.1442645986559986 }
    3 = @{ "G" = 8.542726333333333; "H" = 25.628179; "I" = 0.2146499313812649; "J" = 0.2146499313812649; "O" = 0.007736938980150513; "P" = 0.007736938980150512; "Q" = 9.946461429275777; "R" = 89.51815286348199; "S" = 0.001660733421190341; "T" = 0.001660733421190341 }
    4 = @{ "G" = 8.542726333333333; "H" = 25.628179; "I" = 0.2146499313812649; "J" = 0.2146499313812649; "M" = 26.91044433333333; "N" = 80.73133300000001; "O" = 0.1788207579957193; "P" = 0.1788207579957193; "Q" = 229.8885614480674; "R" = 2068.997053032607; "S" = 0.03838386343332691; "T" = 0.03838386343332691 }
    5 = @{ "G" = 8.542726333333333; "H" = 25.628179; "I" = 0.2146499313812649; "J" = 0.2146499313812649; "M" = 0.7853516666666667; "N" = 2.356055; "O" = 0.00521868678892747; "P" = 0.00521868678892747; "Q" = 6.709044363760555; "R" = 60.38139927384499; "S" = 0.001120190761143595; "T" = 0.001120190761143595 }
    6 = @{ "G" = 8.542726333333333; "H" = 25.628179; "I" = 0.2146499313812649; "J" = 0.2146499313812649; "M" = 20.486157; "N" = 61.458471; "O" = 0.1361311644572737; "P" = 0.1361311644572737; "Q" = 175.007632872701; "R" = 1575.068695854309; "S" = 0.0292205451096055; "T" = 0.02922054510960549 }
    7 = @{ "I" = 0.2649602054889376; "J" = 0.2649602054889376; "M" = 101.1420973333333; "N" = 303.426292; "O" = 0.6720924517779291; "P" = 0.6720924517779291; "Q" = 1066.542944383546; "R" = 9598.886499451912; "S" = 0.178077754130644; "T" = 0.178077754130644 }
    8 = @{ "I" = 0.2649602054889376; "J" = 0.2649602054889376; "O" = 0.007736938980150513; "P" = 0.007736938980150512; "S" = 0.002049980942036051; "T" = 0.002049980942036051 }
    9 = @{ "I" = 0.2649602054889376; "J" = 0.2649602054889376; "M" = 26.91044433333333; "N" = 80.73133300000001; "O" = 0.1788207579957193; "P" = 0.1788207579957193; "Q" = 283.7705099129265; "R" = 2553.934589216338; "S" = 0.04738038478423336; "T" = 0.04738038478423336 }
    10 = @{ "I" = 0.2649602054889376; "J" = 0.2649602054889376; "M" = 0.7853516666666667; "N" = 2.356055; "O" = 0.00521868678892747; "P" = 0.00521868678892747; "Q" = 8.281529660025555; "R" = 74.53376694023; "S" = 0.001382744323976627; "T" = 0.001382744323976626 }
    11 = @{ "I" = 0.2649602054889376; "J" = 0.2649602054889376; "M" = 20.486157; "N" = 61.458471; "O" = 0.1361311644572737; "P" = 0.1361311644572737; "Q" = 216.026429962934; "R" = 1944.237869666406; "S" = 0.03606934130804761; "T" = 0.03606934130804761 }
    12 = @{ "G" = 10.06002866666667; "H" = 30.180086; "I" = 0.2527746270611218; "J" = 0.2527746270611218; "M" = 101.1420973333333; "N" = 303.426292; "O" = 0.6720924517779291; "P" = 0.6720924517779291; "Q" = 1017.492398580123; "R" = 9157.431587221112; "S" = 0.169887918848761; "T" = 0.169887918848761 }
    13 = @{ "G" = 10.06002866666667; "H" = 30.180086; "I" = 0.2527746270611218; "J" = 0.2527746270611218; "O" = 0.007736938980150513; "P" = 0.007736938980150512; "Q" = 11.71308587048755; "R" = 105.417772834388; "S" = 0.001955701865302202; "T" = 0.001955701865302202 }
    14 = @{ "G" = 10.06002866666667; "H" = 30.180086; "I" = 0.2527746270611218; "J" = 0.2527746270611218; "M" = 26.91044433333333; "N" = 80.73133300000001; "O" = 0.1788207579957193; "P" = 0.1788207579957193; "Q" = 270.7198414260709; "R" = 2436.478572834638; "S" = 0.04520135041315505; "T" = 0.04520135041315505 }
    15 = @{ "G" = 10.06002866666667; "H" = 30.180086; "I" = 0.2527746270611218; "J" = 0.2527746270611218; "M" = 0.7853516666666667; "N" = 2.356055; "O" = 0.00521868678892747; "P" = 0.00521868678892747; "Q" = 7.900660280081111; "R" = 71.10594252073; "S" = 0.001319151606819945; "T" = 0.001319151606819944 }
    16 = @{ "G" = 10.06002866666667; "H" = 30.180086; "I" = 0.2527746270611218; "J" = 0.2527746270611218; "M" = 20.486157; "N" = 61.458471; "O" = 0.1361311644572737; "P" = 0.1361311644572737; "Q" = 206.091326689834; "R" = 1854.821940208506; "S" = 0.03441050432708361; "T" = 0.0344105043270836 }
    17 = @{ "G" = 3.225032333333334; "H" = 9.675097000000001; "I" = 0.08103419705149875; "J" = 0.08103419705149875; "M" = 101.1420973333333; "N" = 303.426292; "O" = 0.6720924517779291; "P" = 0.6720924517779291; "Q" = 326.1865341611471; "R" = 2935.678807450324; "S" = 0.05446247217419763; "T" = 0.05446247217419763 }
    18 = @{ "G" = 3.225032333333334; "H" = 9.675097000000001; "I" = 0.08103419705149875; "J" = 0.08103419705149875; "O" = 0.007736938980150513; "P" = 0.007736938980150512; "Q" = 3.754967496325111; "R" = 33.794707466926; "S" = 0.0006269566378929385; "T" = 0.0006269566378929385 }
    19 = @{ "G" = 3.225032333333334; "H" = 9.675097000000001; "I" = 0.08103419705149875; "J" = 0.08103419705149875; "M" = 26.91044433333333; "N" = 80.73133300000001; "O" = 0.1788207579957193; "P" = 0.1788207579957193; "Q" = 86.78705307936679; "R" = 781.0834777143011; "S" = 0.01449059654032349; "T" = 0.01449059654032349 }
    20 = @{ "G" = 3.225032333333334; "H" = 9.675097000000001; "I" = 0.08103419705149875; "J" = 0.08103419705149875; "M" = 0.7853516666666667; "N" = 2.356055; "O" = 0.00521868678892747; "P" = 0.00521868678892747; "Q" = 2.532784518037222; "R" = 22.795060662335; "S" = 0.0004228920936040019; "T" = 0.0004228920936040018 }
    21 = @{ "G" = 3.225032333333334; "H" = 9.675097000000001; "I" = 0.08103419705149875; "J" = 0.08103419705149875; "M" = 20.486157; "N" = 61.458471; "O" = 0.1361311644572737; "P" = 0.1361311644572737; "Q" = 66.06851871074301; "R" = 594.616668396687; "S" = 0.0110312796054807; "T" = 0.0110312796054807 }
    22 = @{ "G" = 7.425629000000001; "H" = 22.276887; "I" = 0.1865810390171769; "J" = 0.1865810390171769; "M" = 101.1420973333333; "N" = 303.426292; "O" = 0.6720924517779291; "P" = 0.6720924517779291; "Q" = 751.0436910792226; "R" = 6759.393219713004; "S" = 0.1253997079683279; "T" = 0.1253997079683279 }
    23 = @{ "G" = 7.425629000000001; "H" = 22.276887; "I" = 0.1865810390171769; "J" = 0.1865810390171769; "O" = 0.007736938980150513; "P" = 0.007736938980150512; "Q" = 8.645803406860667; "R" = 77.812230661746; "S" = 0.00144356611372898; "T" = 0.001443566113728979 }
    24 = @{ "G" = 7.425629000000001; "H" = 22.276887; "I" = 0.1865810390171769; "J" = 0.1865810390171769; "M" = 26.91044433333333; "N" = 80.73133300000001; "O" = 0.1788207579957193; "P" = 0.1788207579957193; "Q" = 199.8269758444857; "R" = 1798.442782600371; "S" = 0.03336456282468044; "T" = 0.03336456282468044 }
    25 = @{ "G" = 7.425629000000001; "H" = 22.276887; "I" = 0.1865810390171769; "J" = 0.1865810390171769; "M" = 0.7853516666666667; "N" = 2.356055; "O" = 0.00521868678892747; "P" = 0.00521868678892747; "Q" = 5.831730111198334; "R" = 52.48557100078501; "S" = 0.0009737080033833019; "T" = 0.0009737080033833017 }
    26 = @{ "G" = 7.425629000000001; "H" = 22.276887; "I" = 0.1865810390171769; "J" = 0.1865810390171769; "M" = 20.486157; "N" = 61.458471; "O" = 0.1361311644572737; "P" = 0.1361311644572737; "Q" = 152.122601517753; "R" = 1369.103413659777; "S" = 0.02539949410705631; "T" = 0.02539949410705631 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
